$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.525.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.444.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.900.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.408.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.445.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "319.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "640.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.574.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0949"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E35").Value = "  -4.29%  "
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "151.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.362"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0307"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "151.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.601"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0501"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0902"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.99%  "
